$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: clear the "Preferred Callback Date/Time/Period" cells (G2, H2, J2) ---
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = ""
$ws.Range("J2").Value = ""

# --- Row 3: new reschedule request for Aarav Mehta ---
$ws.Range("A3").Value = "Aarav Mehta"
$ws.Range("B3").Value = "'917823844614"
$ws.Range("C3").Value = "24 MG Road, Bengaluru"
$ws.Range("D3").Value = "'28"
$ws.Range("E3").Value = "Male"
$ws.Range("F3").Value = "2025-06-25 19:11:27"
$ws.Range("K3").Value = "Pending Callback"
$ws.Range("L3").Value = "Low"

# --- Row 4: new reschedule request for Aarav Mehta ---
$ws.Range("A4").Value = "Aarav Mehta"
$ws.Range("B4").Value = "'917823844614"
$ws.Range("C4").Value = "24 MG Road, Bengaluru"
$ws.Range("D4").Value = "'28"
$ws.Range("E4").Value = "Male"
$ws.Range("F4").Value = "2025-06-25 19:35:25"
$ws.Range("K4").Value = "Pending Callback"
$ws.Range("L4").Value = "Low"

# --- Row 5: new reschedule request for Aarav Mehta ---
$ws.Range("A5").Value = "Aarav Mehta"
$ws.Range("B5").Value = "'917823844614"
$ws.Range("C5").Value = "24 MG Road, Bengaluru"
$ws.Range("D5").Value = "'28"
$ws.Range("E5").Value = "Male"
$ws.Range("F5").Value = "2025-06-25 19:39:22"
$ws.Range("K5").Value = "Pending Callback"
$ws.Range("L5").Value = "Low"

# --- Row 6: new reschedule request for Vanshika panjwani ---
$ws.Range("A6").Value = "Vanshika panjwani"
$ws.Range("B6").Value = "'917823844614"
$ws.Range("C6").Value = "24 MG Road, Bengaluru"
$ws.Range("D6").Value = "'28"
$ws.Range("E6").Value = "Male"
$ws.Range("F6").Value = "2025-06-25 20:34:22"
# G6/H6/I6/J6 stay present but blank (no callback preference captured yet)
$ws.Range("G6").Font.Bold = $false
$ws.Range("H6").Font.Bold = $false
$ws.Range("I6").Font.Bold = $false
$ws.Range("J6").Font.Bold = $false
$ws.Range("K6").Value = "Pending Callback"
$ws.Range("L6").Value = "Low"
